# Refinements to chl a data work up
# Fix the CAL sample label typo: "082516CAL" -> "08252016CAL" (and the "-A" variant)
# These labels live in rows 13 and 14, column A, of the CHL-A sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHL-A")

$ws.Range("A13").Value = "08252016CAL"
$ws.Range("A14").Value = "08252016CAL-A"
